$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) row values for D, J, K, L, M, P columns
$cols = @("D", "J", "K", "L", "M", "P")
$rows = @(2, 3, 4)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New mapping: row2 <- old row3, row3 <- old row4, row4 <- old row2
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$srcRow][$c]
    }
}
